$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.339.47"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.849.97"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'244.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'0.6197"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").Value = "'1.013"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.07470"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'0.2957"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "'23.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "'0.07745"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.830.45"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").Value = "'0.6766"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "'83.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "'0.000009108"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "'5.913"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").Value = "29.321.33"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "2.084.82"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "'239.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.61%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "'7.202"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "'1.016"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'160.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "'0.1439"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").Value = "'8.548"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "'17.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "'1.508"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'0.05630"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("D31").Value = "'4.162"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'4.130"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").Value = "'1.222"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").Value = "'1.856"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "'0.7462"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").Value = "'1.146"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "'2.665"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").Value = "'2.843"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").Value = "'0.01788"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "1.219.30"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").Value = "'6.484"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "'0.9148"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").Value = "'1.012"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "'101.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "1.989.67"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "'65.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "'0.00000000124"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").Value = "'0.5147"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "'0.4064"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "'9.180"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").Value = "'0.05847"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
